# Apply the "Implementação com o Dash" edits to the OEE workbook:
#   - C2: 0.9381 -> 5.9381
#   - F7: 0.6652 -> 8.6652
#   - sheet view: zoom to 175%, move the active selection to G9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two data values that changed.
$ws.Range("C2").Value = 5.9381
$ws.Range("F7").Value = 8.6652

# Zoom the active window to 175% (mirrors zoomScale/zoomScaleNormal in the sheetView).
$excel.ActiveWindow.Zoom = 175

# Move/collapse the selection onto G9 (was J16).
$ws.Range("G9").Select()
